# PRECISE_ENGINEERING.xlsx -- "project model & work order model export remove"
#
# The sheet's header row is repurposed from a Project/Work-Order export
# layout (Created Date, Customer Name, Customer Code, Project Name, Qty)
# to a Customer-detail export layout (Customer Name, Email Id, Code,
# Contact Person, Phone No., Address, GST No) with two extra columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: insert a new (blank) column at E. This shifts the old
#    column E ("Qty", which carries the distinct bold style used only by
#    that header cell) to F, preserving its exact formatting in place so
#    we don't have to fight Excel's style de-duplication later.
# ---------------------------------------------------------------------
[void]$ws.Range("E1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2) Write the new header text. Order matters here only insofar as it
#    matches how the strings were authored; functionally this just sets
#    each header cell's text.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Customer Name"
$ws.Range("B1").Value = "Email Id"
$ws.Range("C1").Value = "Code"
$ws.Range("E1").Value = "Phone No."
$ws.Range("F1").Value = "Address"
$ws.Range("D1").Value = "Contact Person"
$ws.Range("G1").Value = "GST No"

# ---------------------------------------------------------------------
# 3) Formatting.
#    - A1:E1 and G1 are bold headers (same style already used by
#      A1:D1); F1 already carries the correct bold style inherited from
#      the old "Qty" header cell via the column insert above.
#    - Row 2 (the blank data row) uses the plain numeric-format style
#      already present on C2; apply it across the whole row.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("C2").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Column widths for the new A:G layout.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.8333333333333
$ws.Columns.Item(2).ColumnWidth = 33
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 12.6666666666667
$ws.Columns.Item(5).ColumnWidth = 13
$ws.Columns.Item(6).ColumnWidth = 12.8333333333333
$ws.Columns.Item(7).ColumnWidth = 16

# ---------------------------------------------------------------------
# 5) Leave the cursor where the author left it when they saved.
# ---------------------------------------------------------------------
[void]$ws.Range("D11").Select()
